$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataDictionary")

# --- Insert a new row at 36 (pushes current row 36 "MOC Type" etc. down to 37) ---
$ws.Rows.Item(36).Insert()

# --- Populate the new "Profession" row ---
$ws.Range("A36").Value = "Profession"
$ws.Range("B36").Value = "A grouping of Authorised Persons and Disciplines, alllows definition of who is responsible for what discipline scopes."
$ws.Range("C36").Value = "Electrical Technician, Piping Engineer"
$ws.Range("D36").Value = "Reference Tables"
$ws.Range("E36").Value = "C"
$ws.Range("F36").Value = 33
$ws.Range("G36").Value = "Yes"
$ws.Range("H36").Value = "Authorised Person`nDiscipline"
$ws.Range("I36").Value = "As Built Drawing`nCable`nCertification Grouping`nITR`nJob Card`nHandover`nLine`nMechanical Joint`nMOC`nProcedure`nPunch List Item`nPWL`nSpool`nTag`nWalkdown`nWork Pack"

# Formatting for the new row to mirror the rest of the table
$ws.Range("B36:C36").WrapText = $true
$ws.Range("H36").WrapText = $true
$ws.Range("H36").NumberFormat = "d-mmm-yy"
$ws.Rows.Item(36).RowHeight = 264

# --- Expand the table to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:I67"))

# --- Expand the conditional formatting range that flags "Yes" values ---
$cf = $ws.Range("G4:G66").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("G4:G67"))

# --- Fix up the defined names that point past the insertion point ---
$ws.Names.Item("_ftn1").RefersTo = "=DataDictionary!`$A`$95"
$ws.Names.Item("_ftnref1").RefersTo = "=DataDictionary!`$B`$78"
$ws.Names.Item("_Hlk512870131").RefersTo = "=DataDictionary!`$A`$59"

# --- Independent data corrections ---
$ws.Range("F6").Value = 2
$ws.Rows.Item(20).RowHeight = 264

# --- Restore view state ---
$ws.Application.ActiveWindow.ScrollRow = 62
$ws.Range("F29:F67").Select()
